$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "sCs"
$ws.Range("B2").Value = "ECs"
$ws.Range("C2").Value = "FAPs"
$ws.Range("D2").Value = "sCs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.009849666666666666
$ws.Range("H2").Value = 0.029549
$ws.Range("I2").Value = 0.03297364251121477
$ws.Range("J2").Value = 0.03297364251121477
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.02986833333333333
$ws.Range("N2").Value = 0.089605
$ws.Range("O2").Value = 0.02366079079139609
$ws.Range("P2").Value = 0.02366079079139609
$ws.Range("Q2").Value = 0.0002941931272222222
$ws.Range("R2").Value = 0.002647738145
$ws.Range("S2").Value = 0.000780182457088137
$ws.Range("T2").Value = 0.000780182457088137
$ws.Range("A3").Value = "sCs"
$ws.Range("B3").Value = "ECs"
$ws.Range("C3").Value = "FAPs"
$ws.Range("D3").Value = "Cntn2"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.009849666666666666
$ws.Range("H3").Value = 0.029549
$ws.Range("I3").Value = 0.03297364251121477
$ws.Range("J3").Value = 0.03297364251121477
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.1287113333333333
$ws.Range("N3").Value = 0.386134
$ws.Range("O3").Value = 0.1019612275145911
$ws.Range("P3").Value = 0.1019612275145911
$ws.Range("Q3").Value = 0.001267763729555555
$ws.Range("R3").Value = 0.011409873566
$ws.Range("S3").Value = 0.003362033066070763
$ws.Range("T3").Value = 0.003362033066070763
$ws.Range("A4").Value = "sCs"
$ws.Range("B4").Value = "ECs"
$ws.Range("C4").Value = "FAPs"
$ws.Range("D4").Value = "Nrcam"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.009849666666666666
$ws.Range("H4").Value = 0.029549
$ws.Range("I4").Value = 0.03297364251121477
$ws.Range("J4").Value = 0.03297364251121477
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.103776
$ws.Range("N4").Value = 3.311328
$ws.Range("O4").Value = 0.8743779816940128
$ws.Range("P4").Value = 0.8743779816940128
$ws.Range("Q4").Value = 0.01087182567466666
$ws.Range("R4").Value = 0.09784643107199999
$ws.Range("S4").Value = 0.02883142698805587
$ws.Range("T4").Value = 0.02883142698805587
$ws.Range("A5").Value = "Cntn2"
$ws.Range("B5").Value = "ECs"
$ws.Range("C5").Value = "FAPs"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1092446666666667
$ws.Range("H5").Value = 0.327734
$ws.Range("I5").Value = 0.3657174102260808
$ws.Range("J5").Value = 0.3657174102260808
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.02986833333333333
$ws.Range("N5").Value = 0.089605
$ws.Range("O5").Value = 0.02366079079139609
$ws.Range("P5").Value = 0.02366079079139609
$ws.Range("Q5").Value = 0.003262956118888889
$ws.Range("R5").Value = 0.02936660507
$ws.Range("S5").Value = 0.008653163132130479
$ws.Range("T5").Value = 0.008653163132130479
$ws.Range("A6").Value = "Cntn2"
$ws.Range("B6").Value = "ECs"
$ws.Range("C6").Value = "FAPs"
$ws.Range("D6").Value = "Cntn2"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.1092446666666667
$ws.Range("H6").Value = 0.327734
$ws.Range("I6").Value = 0.3657174102260808
$ws.Range("J6").Value = 0.3657174102260808
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1287113333333333
$ws.Range("N6").Value = 0.386134
$ws.Range("O6").Value = 0.1019612275145911
$ws.Range("P6").Value = 0.1019612275145911
$ws.Range("Q6").Value = 0.01406102670622222
$ws.Range("R6").Value = 0.126549240356
$ws.Range("S6").Value = 0.03728899607010848
$ws.Range("T6").Value = 0.03728899607010847
$ws.Range("A7").Value = "Cntn2"
$ws.Range("B7").Value = "ECs"
$ws.Range("C7").Value = "FAPs"
$ws.Range("D7").Value = "Nrcam"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.1092446666666667
$ws.Range("H7").Value = 0.327734
$ws.Range("I7").Value = 0.3657174102260808
$ws.Range("J7").Value = 0.3657174102260808
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.103776
$ws.Range("N7").Value = 3.311328
$ws.Range("O7").Value = 0.8743779816940128
$ws.Range("P7").Value = 0.8743779816940128
$ws.Range("Q7").Value = 0.1205816411946667
$ws.Range("R7").Value = 1.085234770752
$ws.Range("S7").Value = 0.3197752510238419
$ws.Range("T7").Value = 0.3197752510238419
$ws.Range("A8").Value = "Nrcam"
$ws.Range("B8").Value = "ECs"
$ws.Range("C8").Value = "FAPs"
$ws.Range("D8").Value = "sCs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.179619
$ws.Range("H8").Value = 0.5388569999999999
$ws.Range("I8").Value = 0.6013089472627046
$ws.Range("J8").Value = 0.6013089472627045
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.02986833333333333
$ws.Range("N8").Value = 0.089605
$ws.Range("O8").Value = 0.02366079079139609
$ws.Range("P8").Value = 0.02366079079139609
$ws.Range("Q8").Value = 0.005364920164999999
$ws.Range("R8").Value = 0.048284281485
$ws.Range("S8").Value = 0.01422744520217748
$ws.Range("T8").Value = 0.01422744520217747
$ws.Range("A9").Value = "Nrcam"
$ws.Range("B9").Value = "ECs"
$ws.Range("C9").Value = "FAPs"
$ws.Range("D9").Value = "Cntn2"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.179619
$ws.Range("H9").Value = 0.5388569999999999
$ws.Range("I9").Value = 0.6013089472627046
$ws.Range("J9").Value = 0.6013089472627045
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.1287113333333333
$ws.Range("N9").Value = 0.386134
$ws.Range("O9").Value = 0.1019612275145911
$ws.Range("P9").Value = 0.1019612275145911
$ws.Range("Q9").Value = 0.02311900098199999
$ws.Range("R9").Value = 0.2080710088379999
$ws.Range("S9").Value = 0.06131019837841189
$ws.Range("T9").Value = 0.06131019837841187
$ws.Range("A10").Value = "Nrcam"
$ws.Range("B10").Value = "ECs"
$ws.Range("C10").Value = "FAPs"
$ws.Range("D10").Value = "Nrcam"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.179619
$ws.Range("H10").Value = 0.5388569999999999
$ws.Range("I10").Value = 0.6013089472627046
$ws.Range("J10").Value = 0.6013089472627045
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.103776
$ws.Range("N10").Value = 3.311328
$ws.Range("O10").Value = 0.8743779816940128
$ws.Range("P10").Value = 0.8743779816940128
$ws.Range("Q10").Value = 0.1982591413439999
$ws.Range("R10").Value = 1.784332272096
$ws.Range("S10").Value = 0.5257713036821152
$ws.Range("T10").Value = 0.5257713036821151
